# Update the cryptocurrency price/volume snapshot (and the two swapped
# ranking rows) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text values (e.g. "29.045.52") in the
# source workbook. A bare numeric-looking string assigned via .Value
# would be auto-converted to a real number by Excel, so the whole
# column is switched to Text format for the duration of the writes and
# then restored to the default style afterwards.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.999.71"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "1.898.91"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "326.99"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "0.4612"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("D9").Value = "0.07883"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "0.9938"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("D11").Value = "22.14"
$ws.Range("E11").Value = "  +2.67%  "

$ws.Range("D12").Value = "1.856.49"
$ws.Range("E12").Value = "  -5.01%  "

$ws.Range("D13").Value = "7.069"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("D14").Value = "5.726"
$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("D15").Value = "0.06962"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "88.34"
$ws.Range("E16").Value = "  +1.79%  "

$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "0.00001001"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "16.99"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").Value = "29.007.19"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("D22").Value = "5.328"
$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  +0.85%  "

$ws.Range("D24").Value = "2.160.77"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "2.067"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").Value = "156.35"
$ws.Range("E26").Value = "  +2.01%  "

$ws.Range("D27").Value = "19.37"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").Value = "5.986"
$ws.Range("E28").Value = "  +5.76%  "

$ws.Range("D29").Value = "1.934"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").Value = "118.00"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "0.09349"
$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("D32").Value = "0.9142"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").Value = "5.326"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").Value = "1.342"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").Value = "3.289"
$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("E36").Value = "  +4.86%  "

$ws.Range("D37").Value = "0.05791"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("D38").Value = "0.02086"
$ws.Range("E38").Value = "  +1.44%  "

$ws.Range("B39").Value = "Frax"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D39").Value = "1.003"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.824"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").Value = "0.5703"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").Value = "0.1782"
$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").Value = "9.848"
$ws.Range("E43").Value = "  +1.89%  "

$ws.Range("D44").Value = "2.268"
$ws.Range("E44").Value = "  +5.57%  "

$ws.Range("D45").Value = "12.01"
$ws.Range("E45").Value = "  +3.36%  "

$ws.Range("D46").Value = "0.5367"
$ws.Range("E46").Value = "  +2.46%  "

$ws.Range("D47").Value = "0.07046"
$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("D48").Value = "1.851"
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "112.88"
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.533"
$ws.Range("E50").Value = "  +5.27%  "

$ws.Range("D51").Value = "1.071"
$ws.Range("E51").Value = "  -4.25%  "

# Restore the default (unformatted) style on column D so only the
# values themselves changed.
$ws.Range("D2:D51").Style = "Normal"
